$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I2").NumberFormat = "mm/dd/yy"
Write-Host ("NumberFormat I2: " + $ws.Range("I2").NumberFormat)
